$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row22 = @(0.21, 4693.743604618, 34.74120455, 4659.002400068, 25.32429377533333, 1764.728998344, 21.748509934, 1742.98048841, 26.59494055133333)
for ($c = 1; $c -le 9; $c++) { $ws.Cells.Item(22, $c).Value = $row22[$c-1] }

$row23 = @(0.22, 4873.862068222, 34.549828712, 4839.31223951, 25.45483019533333, 1814.294152412, 21.730126024, 1792.564026386, 26.73886206266667)
for ($c = 1; $c -le 9; $c++) { $ws.Cells.Item(23, $c).Value = $row23[$c-1] }

$row24 = @(0.23, 4989.959616158, 34.193174196, 4955.766441962, 25.57009584133333, 1834.75095213, 21.69227645, 1813.058675678, 26.86194370533333)
for ($c = 1; $c -le 9; $c++) { $ws.Cells.Item(24, $c).Value = $row24[$c-1] }

$row25 = @(0.24, 4807.966401793999, 33.930818074, 4774.035583722, 25.81990274533334, 1745.358372966, 21.805460372, 1723.552912592, 27.12114217266666)
for ($c = 1; $c -le 9; $c++) { $ws.Cells.Item(25, $c).Value = $row25[$c-1] }

$row26 = @(0.25, 4828.259401284, 33.617331686, 4794.642069598, 25.91115374533333, 1739.354431102, 21.88837518, 1717.466055924, 27.21708276466667)
for ($c = 1; $c -le 9; $c++) { $ws.Cells.Item(26, $c).Value = $row26[$c-1] }

$row27 = @(0.26, 4866.381600115999, 33.774286282, 4832.607313834001, 26.42515303066667, 1736.166937038, 22.310224142, 1713.856712894, 27.74777298466667)
for ($c = 1; $c -le 9; $c++) { $ws.Cells.Item(27, $c).Value = $row27[$c-1] }

$row28 = @(0.27, 4933.811597968, 33.967502006, 4899.844095962, 27.12920759, 1742.295888852, 23.002376966, 1719.293511884, 28.46061579466667)
for ($c = 1; $c -le 9; $c++) { $ws.Cells.Item(28, $c).Value = $row28[$c-1] }
